$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new data points for TimSortDLL (column D) for the 10000 and 100000 rows
$ws.Range("D4").Value = 0.088007450103759696
$ws.Range("D5").Value = 1.1700975894927901

# Update the active selection as it was left in the original workbook
$ws.Range("E11").Select()
